# Weekly update: insert a fresh week's record (row 139) for the
# "Hortaliza, Terminal La Palmera de La Serena - Jengibre" sheet.
# Inserting the row shifts every following record down by one,
# which reproduces the observed diff (old row 139 -> new row 140,
# ..., old row 195 -> new row 196) without touching rows 1-138.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all data rows from 139 downward by one row.
$ws.Rows.Item(139).Insert()

# Populate the newly inserted row with this week's observation.
$ws.Cells.Item(139, 1).Value2  = 8
$ws.Cells.Item(139, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(139, 3).Value2  = "Coquimbo"
$ws.Cells.Item(139, 4).Value2  = 45229
$ws.Cells.Item(139, 5).Value2  = 4
$ws.Cells.Item(139, 6).Value2  = 100114007
$ws.Cells.Item(139, 7).Value2  = "Jengibre"
$ws.Cells.Item(139, 8).Value2  = "Sin especificar"
$ws.Cells.Item(139, 9).Value2  = "Primera"
$ws.Cells.Item(139, 10).Value2 = 500
$ws.Cells.Item(139, 11).Value2 = 27000
$ws.Cells.Item(139, 12).Value2 = 28000
$ws.Cells.Item(139, 13).Value2 = 27500
$ws.Cells.Item(139, 14).Value2 = "`$/caja 13 kilos"
$ws.Cells.Item(139, 15).Value2 = "Perú"
$ws.Cells.Item(139, 16).Value2 = 2115
$ws.Cells.Item(139, 17).Value2 = 13
$ws.Cells.Item(139, 18).Value2 = "Hortaliza"
